$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 14707548
$ws.Range("I33").Value = 25000674
$ws.Range("J33").Value = 3084.2856
$ws.Range("K33").Value = 25000674
$ws.Range("L33").Value = 3084.2856
$ws.Range("M33").Value = -25000445
$ws.Range("N33").Value = -3542.2856

$ws.Range("H55").Value = 253.21053
$ws.Range("I55").Value = 260.57144
$ws.Range("J55").Value = 248.91667
$ws.Range("K55").Value = 260.57144
$ws.Range("L55").Value = 248.91667
$ws.Range("M55").Value = -46.57144
$ws.Range("N55").Value = -676.9166700000001

$ws.Range("H70").Value = 1766.3334
$ws.Range("J70").Value = 1879.6
$ws.Range("L70").Value = 5638.799999999999
$ws.Range("N70").Value = -6178.799999999999

$ws.Range("H73").Value = 1766.3334
$ws.Range("J73").Value = 1879.6
$ws.Range("L73").Value = 5638.799999999999
$ws.Range("N73").Value = -7510.799999999999

$ws.Range("H92").Value = 2654.5557
$ws.Range("I92").Value = 2364
$ws.Range("J92").Value = 3410
$ws.Range("K92").Value = 2364
$ws.Range("L92").Value = 3410
$ws.Range("M92").Value = -1116
$ws.Range("N92").Value = -5906

$ws.Range("H103").Value = 899
$ws.Range("I103").Value = 998
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 2994
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -2408
$ws.Range("N103").Value = -3572

$ws.Range("H132").Value = 17252.578
$ws.Range("I132").Value = 12099.944
$ws.Range("J132").Value = 110000
$ws.Range("K132").Value = 36299.83199999999
$ws.Range("L132").Value = 330000
$ws.Range("M132").Value = -33769.83199999999
$ws.Range("N132").Value = -335060

$ws.Range("H135").Value = 1366.5
$ws.Range("I135").Value = 1362.7778
$ws.Range("K135").Value = 12265.0002
$ws.Range("M135").Value = -9730.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1128767.4
$ws.Range("I61").Value = 3487.8333
$ws.Range("K61").Value = 3487.8333
$ws.Range("M61").Value = -3275.8333

$ws.Range("H132").Value = 3157989.8
$ws.Range("I132").Value = 2018.5454
$ws.Range("K132").Value = 6055.6362
$ws.Range("M132").Value = -3525.6362

$ws.Range("H136").Value = 1128767.4
$ws.Range("I136").Value = 3487.8333
$ws.Range("K136").Value = 10463.4999
$ws.Range("M136").Value = -7913.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 497.5
$ws.Range("I8").Value = 497.5
$ws.Range("K8").Value = 497.5
$ws.Range("M8").Value = -357.5

$ws.Range("H16").Value = 3769.3333
$ws.Range("I16").Value = 4404
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 4404
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -4234
$ws.Range("N16").Value = -2840

$ws.Range("H22").Value = 1296.4445
$ws.Range("I22").Value = 1396
$ws.Range("K22").Value = 1396
$ws.Range("M22").Value = -1223

$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = ""

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1088.7
$ws.Range("I16").Value = 1113.5714
$ws.Range("J16").Value = 1030.6666
$ws.Range("K16").Value = 1113.5714
$ws.Range("L16").Value = 1030.6666
$ws.Range("M16").Value = -826.5714
$ws.Range("N16").Value = -1604.6666

$ws.Range("H31").Value = 55094.242
$ws.Range("I31").Value = 70095.83
$ws.Range("K31").Value = 70095.83
$ws.Range("M31").Value = -69800.83

$ws.Range("H34").Value = 55094.242
$ws.Range("I34").Value = 70095.83
$ws.Range("K34").Value = 70095.83
$ws.Range("M34").Value = -69893.83

$ws.Range("H113").Value = 1088.7
$ws.Range("I113").Value = 1113.5714
$ws.Range("J113").Value = 1030.6666
$ws.Range("K113").Value = 1113.5714
$ws.Range("L113").Value = 1030.6666
$ws.Range("M113").Value = 1056.4286
$ws.Range("N113").Value = -5370.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1343.5
$ws.Range("I129").Value = 1024.0834
$ws.Range("J129").Value = 3260
$ws.Range("K129").Value = 3072.2502
$ws.Range("L129").Value = 9780
$ws.Range("M129").Value = 1927.7498
$ws.Range("N129").Value = -19780

$ws.Range("H131").Value = 1503.433
$ws.Range("I131").Value = 1224.75
$ws.Range("J131").Value = 1515.4193
$ws.Range("K131").Value = 3674.25
$ws.Range("L131").Value = 4546.257900000001
$ws.Range("M131").Value = 1365.75
$ws.Range("N131").Value = -14626.2579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""

$ws.Range("H54").Value = 34999
$ws.Range("J54").Value = 34999
$ws.Range("L54").Value = 34999
$ws.Range("N54").Value = -35779

$ws.Range("H97").Value = 1127.2142
$ws.Range("I97").Value = 1044.6842
$ws.Range("K97").Value = 1044.6842
$ws.Range("M97").Value = -548.6841999999999

$ws.Range("H122").Value = 2715.8096
$ws.Range("I122").Value = 2836.2222
$ws.Range("K122").Value = 8508.6666
$ws.Range("M122").Value = -6058.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3390.0667
$ws.Range("J40").Value = 4888.3335
$ws.Range("L40").Value = 4888.3335
$ws.Range("N40").Value = -5160.3335

$ws.Range("H46").Value = 3058.1765
$ws.Range("I46").Value = 1332.6666
$ws.Range("J46").Value = 3427.9285
$ws.Range("K46").Value = 1332.6666
$ws.Range("L46").Value = 3427.9285
$ws.Range("M46").Value = -1144.6666
$ws.Range("N46").Value = -3803.9285

$ws.Range("H93").Value = 9599.333000000001
$ws.Range("J93").Value = 1449
$ws.Range("L93").Value = 1449
$ws.Range("N93").Value = -3945

$ws.Range("H94").Value = 26875
$ws.Range("J94").Value = 26875
$ws.Range("L94").Value = 26875
$ws.Range("N94").Value = -28227

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""

$ws.Range("H96").Value = 1522.1111
$ws.Range("I96").Value = 1874.5
$ws.Range("K96").Value = 1874.5
$ws.Range("M96").Value = -501.5

$ws.Range("H126").Value = 3608.3
$ws.Range("I126").Value = 4240.4287
$ws.Range("K126").Value = 12721.2861
$ws.Range("M126").Value = -10251.2861

$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

